$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 8, pushing the existing data
# (old rows 8..190) down to rows 10..192. This also grows the sheet's
# used range / dimension from A1:R190 to A1:R192.
$ws.Range("A8:A9").EntireRow.Insert()

# Fill the two newly inserted rows with the new "Cilantro" price entries
# (Vega Monumental Concepción, fecha 2022-05-18, Región de Ñuble).

# Row 8: Calidad "Primera"
$ws.Cells.Item(8, 1).Value2  = 11
$ws.Cells.Item(8, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value2  = "Bíobío"
$ws.Cells.Item(8, 4).Value2  = 44699
$ws.Cells.Item(8, 5).Value2  = 8
$ws.Cells.Item(8, 6).Value2  = 100112040
$ws.Cells.Item(8, 7).Value2  = "Cilantro"
$ws.Cells.Item(8, 8).Value2  = "Sin especificar"
$ws.Cells.Item(8, 9).Value2  = "Primera"
$ws.Cells.Item(8, 10).Value2 = 200
$ws.Cells.Item(8, 11).Value2 = 600
$ws.Cells.Item(8, 12).Value2 = 700
$ws.Cells.Item(8, 13).Value2 = 650
$ws.Cells.Item(8, 14).Value2 = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(8, 16).Value2 = 650
$ws.Cells.Item(8, 17).Value2 = 1
$ws.Cells.Item(8, 18).Value2 = "Hortaliza"

# Row 9: Calidad "Segunda"
$ws.Cells.Item(9, 1).Value2  = 11
$ws.Cells.Item(9, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value2  = "Bíobío"
$ws.Cells.Item(9, 4).Value2  = 44699
$ws.Cells.Item(9, 5).Value2  = 8
$ws.Cells.Item(9, 6).Value2  = 100112040
$ws.Cells.Item(9, 7).Value2  = "Cilantro"
$ws.Cells.Item(9, 8).Value2  = "Sin especificar"
$ws.Cells.Item(9, 9).Value2  = "Segunda"
$ws.Cells.Item(9, 10).Value2 = 100
$ws.Cells.Item(9, 11).Value2 = 500
$ws.Cells.Item(9, 12).Value2 = 500
$ws.Cells.Item(9, 13).Value2 = 500
$ws.Cells.Item(9, 14).Value2 = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(9, 16).Value2 = 500
$ws.Cells.Item(9, 17).Value2 = 1
$ws.Cells.Item(9, 18).Value2 = "Hortaliza"
